$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.332.93"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.877.57"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'0.7105"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'242.27"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.07995"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("D9").Value = "'0.3155"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("D10").Value = "'24.94"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").Value = "'0.08283"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "1.899.14"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "'5.248"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "'94.38"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("D15").Value = "'0.7127"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "'6.364"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").Value = "'0.000008513"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "29.356.22"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'244.10"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").Value = "2.140.95"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "'7.776"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("E25").Value = "  -2.16%  "
$ws.Range("D26").Value = "'9.065"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'162.63"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'1.504"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").Value = "'4.414"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'4.323"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -7.90%  "
$ws.Range("D33").Value = "'0.05373"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").Value = "'1.933"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "'0.7641"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").Value = "'1.183"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").Value = "'0.01882"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "1.259.11"
$ws.Range("E39").Value = "  +3.38%  "
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "'6.512"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").Value = "'112.76"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Value = "'0.9052"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "'74.18"
$ws.Range("E45").Value = "  +8.39%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "2.029.20"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'0.5227"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").Value = "'1.798"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'9.447"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'0.4371"
$ws.Range("E51").Value = "  +1.22%  "
